# Highlight (green) a block of "Les Rooms" / "Les Comptes" backlog bullet
# points in green, matching the author's commit "Chat presque terminé,
# test a la maison". Both the run text and the paragraph mark (so the
# whole line, including the bullet, renders highlighted) receive the
# wdGreen highlight color.

$d = $word.ActiveDocument

# WdColorIndex.wdGreen
$wdGreen = 4

$targets = @(
    "Les Rooms ",
    "Modifier le serveur pour qu'il gère plusieurs grilles en mémoire ",
    "Faire marcher la commande socket.join(roomId) : séparer les joueurs dans des canaux différents.",
    "La Liste Publique",
    "Afficher sur l'accueil la liste des grilles créées par les autres.",
    "Ajouter le bouton `"Rejoindre`" qui connecte le joueur à la bonne Room.",
    "Les Comptes (Base de données)",
    "Installer et connecter MongoDB.",
    "Créer le formulaire d'inscription / connexion simple (Pseudo + Mot de passe)."
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    foreach ($target in $targets) {
        if ($t -eq ($target + [char]13)) {
            # Setting Font.HighlightColorIndex on the paragraph's own Range
            # colors the run(s) AND folds the highlight onto the paragraph
            # mark's run properties (w:pPr/w:rPr), matching how Word stores
            # a highlighted bullet line.
            $p.Range.Font.HighlightColorIndex = $wdGreen
        }
    }
}
